$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.380719000000001
$ws.Range("H2").Value = 28.142157
$ws.Range("I2").Value = 0.03679977590837273
$ws.Range("J2").Value = 0.03679977590837273
$ws.Range("M2").Value = 0.177232
$ws.Range("N2").Value = 0.531696
$ws.Range("O2").Value = 0.0005104719838156216
$ws.Range("P2").Value = 0.0005104719838156217
$ws.Range("Q2").Value = 1.662563589808
$ws.Range("R2").Value = 14.963072308272
$ws.Range("S2").Value = 0.00001878525461191735
$ws.Range("T2").Value = 0.00001878525461191735
$ws.Range("G3").Value = 9.380719000000001
$ws.Range("H3").Value = 28.142157
$ws.Range("I3").Value = 0.03679977590837273
$ws.Range("J3").Value = 0.03679977590837273
$ws.Range("M3").Value = 0.08113566666666668
$ws.Range("N3").Value = 0.243407
$ws.Range("O3").Value = 0.0002336907822601807
$ws.Range("P3").Value = 0.0002336907822601807
$ws.Range("Q3").Value = 0.7611108898776668
$ws.Range("R3").Value = 6.849998008899001
$ws.Range("S3").Value = 0.000008599768419026976
$ws.Range("T3").Value = 0.000008599768419026976
$ws.Range("G4").Value = 9.380719000000001
$ws.Range("H4").Value = 28.142157
$ws.Range("I4").Value = 0.03679977590837273
$ws.Range("J4").Value = 0.03679977590837273
$ws.Range("M4").Value = 274.5137023333334
$ws.Range("N4").Value = 823.541107
$ws.Range("O4").Value = 0.7906673411949746
$ws.Range("P4").Value = 0.7906673411949746
$ws.Range("Q4").Value = 2575.135903238645
$ws.Range("R4").Value = 23176.2231291478
$ws.Range("S4").Value = 0.02909638097404395
$ws.Range("T4").Value = 0.02909638097404395
$ws.Range("G5").Value = 9.380719000000001
$ws.Range("H5").Value = 28.142157
$ws.Range("I5").Value = 0.03679977590837273
$ws.Range("J5").Value = 0.03679977590837273
$ws.Range("M5").Value = 0.042388
$ws.Range("N5").Value = 0.127164
$ws.Range("O5").Value = 0.0001220879211991998
$ws.Range("P5").Value = 0.0001220879211991998
$ws.Range("Q5").Value = 0.3976299169720001
$ws.Range("R5").Value = 3.578669252748
$ws.Range("S5").Value = 0.00000449280814124962
$ws.Range("T5").Value = 0.00000449280814124962
$ws.Range("G6").Value = 9.380719000000001
$ws.Range("H6").Value = 28.142157
$ws.Range("I6").Value = 0.03679977590837273
$ws.Range("J6").Value = 0.03679977590837273
$ws.Range("M6").Value = 72.37795533333333
$ws.Range("N6").Value = 217.133866
$ws.Range("O6").Value = 0.2084664081177503
$ws.Range("P6").Value = 0.2084664081177503
$ws.Range("Q6").Value = 678.9572607765514
$ws.Range("R6").Value = 6110.615346988962
$ws.Range("S6").Value = 0.007671517103156586
$ws.Range("T6").Value = 0.007671517103156586
$ws.Range("G7").Value = 90.439374
$ws.Range("H7").Value = 271.318122
$ws.Range("I7").Value = 0.3547860986448385
$ws.Range("J7").Value = 0.3547860986448385
$ws.Range("M7").Value = 0.177232
$ws.Range("N7").Value = 0.531696
$ws.Range("O7").Value = 0.0005104719838156216
$ws.Range("P7").Value = 0.0005104719838156217
$ws.Range("Q7").Value = 16.028751132768
$ws.Range("R7").Value = 144.258760194912
$ws.Range("S7").Value = 0.0001811083636054355
$ws.Range("T7").Value = 0.0001811083636054356
$ws.Range("G8").Value = 90.439374
$ws.Range("H8").Value = 271.318122
$ws.Range("I8").Value = 0.3547860986448385
$ws.Range("J8").Value = 0.3547860986448385
$ws.Range("M8").Value = 0.08113566666666668
$ws.Range("N8").Value = 0.243407
$ws.Range("O8").Value = 0.0002336907822601807
$ws.Range("P8").Value = 0.0002336907822601807
$ws.Range("Q8").Value = 7.337858902406001
$ws.Range("R8").Value = 66.04073012165401
$ws.Range("S8").Value = 0.00008291024092734996
$ws.Range("T8").Value = 0.00008291024092734996
$ws.Range("G9").Value = 90.439374
$ws.Range("H9").Value = 271.318122
$ws.Range("I9").Value = 0.3547860986448385
$ws.Range("J9").Value = 0.3547860986448385
$ws.Range("M9").Value = 274.5137023333334
$ws.Range("N9").Value = 823.541107
$ws.Range("O9").Value = 0.7906673411949746
$ws.Range("P9").Value = 0.7906673411949746
$ws.Range("Q9").Value = 24826.84739344901
$ws.Range("R9").Value = 223441.6265410411
$ws.Range("S9").Value = 0.2805177813084525
$ws.Range("T9").Value = 0.2805177813084525
$ws.Range("G10").Value = 90.439374
$ws.Range("H10").Value = 271.318122
$ws.Range("I10").Value = 0.3547860986448385
$ws.Range("J10").Value = 0.3547860986448385
$ws.Range("M10").Value = 0.042388
$ws.Range("N10").Value = 0.127164
$ws.Range("O10").Value = 0.0001220879211991998
$ws.Range("P10").Value = 0.0001220879211991998
$ws.Range("Q10").Value = 3.833544185112
$ws.Range("R10").Value = 34.501897666008
$ws.Range("S10").Value = 0.00004331509725392256
$ws.Range("T10").Value = 0.00004331509725392256
$ws.Range("G11").Value = 90.439374
$ws.Range("H11").Value = 271.318122
$ws.Range("I11").Value = 0.3547860986448385
$ws.Range("J11").Value = 0.3547860986448385
$ws.Range("M11").Value = 72.37795533333333
$ws.Range("N11").Value = 217.133866
$ws.Range("O11").Value = 0.2084664081177503
$ws.Range("P11").Value = 0.2084664081177503
$ws.Range("Q11").Value = 6545.816971746628
$ws.Range("R11").Value = 58912.35274571966
$ws.Range("S11").Value = 0.07396098363459933
$ws.Range("T11").Value = 0.07396098363459933
$ws.Range("G12").Value = 100.179423
$ws.Range("H12").Value = 300.538269
$ws.Range("I12").Value = 0.3929954960840508
$ws.Range("J12").Value = 0.3929954960840508
$ws.Range("M12").Value = 0.177232
$ws.Range("N12").Value = 0.531696
$ws.Range("O12").Value = 0.0005104719838156216
$ws.Range("P12").Value = 0.0005104719838156217
$ws.Range("Q12").Value = 17.754999497136
$ws.Range("R12").Value = 159.794995474224
$ws.Range("S12").Value = 0.0002006131905166297
$ws.Range("T12").Value = 0.0002006131905166298
$ws.Range("G13").Value = 100.179423
$ws.Range("H13").Value = 300.538269
$ws.Range("I13").Value = 0.3929954960840508
$ws.Range("J13").Value = 0.3929954960840508
$ws.Range("M13").Value = 0.08113566666666668
$ws.Range("N13").Value = 0.243407
$ws.Range("O13").Value = 0.0002336907822601807
$ws.Range("P13").Value = 0.0002336907822601807
$ws.Range("Q13").Value = 8.128124271387001
$ws.Range("R13").Value = 73.153118442483
$ws.Range("S13").Value = 0.00009183942490460961
$ws.Range("T13").Value = 0.00009183942490460961
$ws.Range("G14").Value = 100.179423
$ws.Range("H14").Value = 300.538269
$ws.Range("I14").Value = 0.3929954960840508
$ws.Range("J14").Value = 0.3929954960840508
$ws.Range("M14").Value = 274.5137023333334
$ws.Range("N14").Value = 823.541107
$ws.Range("O14").Value = 0.7906673411949746
$ws.Range("P14").Value = 0.7906673411949746
$ws.Range("Q14").Value = 27500.62430534709
$ws.Range("R14").Value = 247505.6187481238
$ws.Range("S14").Value = 0.3107287039903765
$ws.Range("T14").Value = 0.3107287039903765
$ws.Range("G15").Value = 100.179423
$ws.Range("H15").Value = 300.538269
$ws.Range("I15").Value = 0.3929954960840508
$ws.Range("J15").Value = 0.3929954960840508
$ws.Range("M15").Value = 0.042388
$ws.Range("N15").Value = 0.127164
$ws.Range("O15").Value = 0.0001220879211991998
$ws.Range("P15").Value = 0.0001220879211991998
$ws.Range("Q15").Value = 4.246405382124
$ws.Range("R15").Value = 38.217648439116
$ws.Range("S15").Value = 0.00004798000315755001
$ws.Range("T15").Value = 0.00004798000315755001
$ws.Range("G16").Value = 100.179423
$ws.Range("H16").Value = 300.538269
$ws.Range("I16").Value = 0.3929954960840508
$ws.Range("J16").Value = 0.3929954960840508
$ws.Range("M16").Value = 72.37795533333333
$ws.Range("N16").Value = 217.133866
$ws.Range("O16").Value = 0.2084664081177503
$ws.Range("P16").Value = 0.2084664081177503
$ws.Range("Q16").Value = 7250.781803213106
$ws.Range("R16").Value = 65257.03622891796
$ws.Range("S16").Value = 0.08192635947509548
$ws.Range("T16").Value = 0.08192635947509548
$ws.Range("G17").Value = 1.427630666666667
$ws.Range("H17").Value = 4.282892
$ws.Range("I17").Value = 0.005600475679236752
$ws.Range("J17").Value = 0.005600475679236752
$ws.Range("M17").Value = 0.177232
$ws.Range("N17").Value = 0.531696
$ws.Range("O17").Value = 0.0005104719838156216
$ws.Range("P17").Value = 0.0005104719838156217
$ws.Range("Q17").Value = 0.2530218383146666
$ws.Range("R17").Value = 2.277196544832
$ws.Range("S17").Value = 0.000002858885930291125
$ws.Range("T17").Value = 0.000002858885930291126
$ws.Range("G18").Value = 1.427630666666667
$ws.Range("H18").Value = 4.282892
$ws.Range("I18").Value = 0.005600475679236752
$ws.Range("J18").Value = 0.005600475679236752
$ws.Range("M18").Value = 0.08113566666666668
$ws.Range("N18").Value = 0.243407
$ws.Range("O18").Value = 0.0002336907822601807
$ws.Range("P18").Value = 0.0002336907822601807
$ws.Range("Q18").Value = 0.1158317658937778
$ws.Range("R18").Value = 1.042485893044
$ws.Range("S18").Value = 0.000001308779542509953
$ws.Range("T18").Value = 0.000001308779542509953
$ws.Range("G19").Value = 1.427630666666667
$ws.Range("H19").Value = 4.282892
$ws.Range("I19").Value = 0.005600475679236752
$ws.Range("J19").Value = 0.005600475679236752
$ws.Range("M19").Value = 274.5137023333334
$ws.Range("N19").Value = 823.541107
$ws.Range("O19").Value = 0.7906673411949746
$ws.Range("P19").Value = 0.7906673411949746
$ws.Range("Q19").Value = 391.9041798712716
$ws.Range("R19").Value = 3527.137618841444
$ws.Range("S19").Value = 0.004428113214729242
$ws.Range("T19").Value = 0.004428113214729242
$ws.Range("G20").Value = 1.427630666666667
$ws.Range("H20").Value = 4.282892
$ws.Range("I20").Value = 0.005600475679236752
$ws.Range("J20").Value = 0.005600475679236752
$ws.Range("M20").Value = 0.042388
$ws.Range("N20").Value = 0.127164
$ws.Range("O20").Value = 0.0001220879211991998
$ws.Range("P20").Value = 0.0001220879211991998
$ws.Range("Q20").Value = 0.06051440869866667
$ws.Range("R20").Value = 0.544629678288
$ws.Range("S20").Value = 0.0000006837504334046913
$ws.Range("T20").Value = 0.0000006837504334046913
$ws.Range("G21").Value = 1.427630666666667
$ws.Range("H21").Value = 4.282892
$ws.Range("I21").Value = 0.005600475679236752
$ws.Range("J21").Value = 0.005600475679236752
$ws.Range("M21").Value = 72.37795533333333
$ws.Range("N21").Value = 217.133866
$ws.Range("O21").Value = 0.2084664081177503
$ws.Range("P21").Value = 0.2084664081177503
$ws.Range("Q21").Value = 103.3289886244969
$ws.Range("R21").Value = 929.9608976204721
$ws.Range("S21").Value = 0.001167511048601304
$ws.Range("T21").Value = 0.001167511048601304
$ws.Range("G22").Value = 53.48524799999999
$ws.Range("H22").Value = 160.455744
$ws.Range("I22").Value = 0.2098181536835013
$ws.Range("J22").Value = 0.2098181536835013
$ws.Range("M22").Value = 0.177232
$ws.Range("N22").Value = 0.531696
$ws.Range("O22").Value = 0.0005104719838156216
$ws.Range("P22").Value = 0.0005104719838156217
$ws.Range("Q22").Value = 9.479297473535997
$ws.Range("R22").Value = 85.31367726182398
$ws.Range("S22").Value = 0.0001071062891513479
$ws.Range("T22").Value = 0.0001071062891513479
$ws.Range("G23").Value = 53.48524799999999
$ws.Range("H23").Value = 160.455744
$ws.Range("I23").Value = 0.2098181536835013
$ws.Range("J23").Value = 0.2098181536835013
$ws.Range("M23").Value = 0.08113566666666668
$ws.Range("N23").Value = 0.243407
$ws.Range("O23").Value = 0.0002336907822601807
$ws.Range("P23").Value = 0.0002336907822601807
$ws.Range("Q23").Value = 4.339561253312
$ws.Range("R23").Value = 39.056051279808
$ws.Range("S23").Value = 0.00004903256846668423
$ws.Range("T23").Value = 0.00004903256846668423
$ws.Range("G24").Value = 53.48524799999999
$ws.Range("H24").Value = 160.455744
$ws.Range("I24").Value = 0.2098181536835013
$ws.Range("J24").Value = 0.2098181536835013
$ws.Range("M24").Value = 274.5137023333334
$ws.Range("N24").Value = 823.541107
$ws.Range("O24").Value = 0.7906673411949746
$ws.Range("P24").Value = 0.7906673411949746
$ws.Range("Q24").Value = 14682.43344869651
$ws.Range("R24").Value = 132141.9010382686
$ws.Range("S24").Value = 0.1658963617073725
$ws.Range("T24").Value = 0.1658963617073725
$ws.Range("G25").Value = 53.48524799999999
$ws.Range("H25").Value = 160.455744
$ws.Range("I25").Value = 0.2098181536835013
$ws.Range("J25").Value = 0.2098181536835013
$ws.Range("M25").Value = 0.042388
$ws.Range("N25").Value = 0.127164
$ws.Range("O25").Value = 0.0001220879211991998
$ws.Range("P25").Value = 0.0001220879211991998
$ws.Range("Q25").Value = 2.267132692224
$ws.Range("R25").Value = 20.404194230016
$ws.Range("S25").Value = 0.00002561626221307289
$ws.Range("T25").Value = 0.00002561626221307289
$ws.Range("G26").Value = 53.48524799999999
$ws.Range("H26").Value = 160.455744
$ws.Range("I26").Value = 0.2098181536835013
$ws.Range("J26").Value = 0.2098181536835013
$ws.Range("M26").Value = 72.37795533333333
$ws.Range("N26").Value = 217.133866
$ws.Range("O26").Value = 0.2084664081177503
$ws.Range("P26").Value = 0.2084664081177503
$ws.Range("Q26").Value = 3871.152890736255
$ws.Range("R26").Value = 34840.3760166263
$ws.Range("S26").Value = 0.04374003685629763
$ws.Range("T26").Value = 0.04374003685629763
